$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sample/protocol id shown in column G (shared string) ---
$ws.Range("G2:G37").Value = "E7420"

# --- Re-style column G (was style 1 / general-aligned Arial11) to a new
#     explicit style: Arial 11, black, general alignment ---
$gRange = $ws.Range("G2:G37")
$gRange.Font.Name = "Arial"
$gRange.Font.Size = 11
$gRange.Font.Color = 0
$gRange.HorizontalAlignment = 1
$gRange.NumberFormat = "General"

# --- Column H: turn the stored boolean literal into a live formula
#     "=FALSE()"; keep its existing font/number format untouched ---
$hRange = $ws.Range("H2:H37")
$hRange.Formula = "=FALSE()"

# --- Selection / scroll position used to sit on H2:H37 viewing row 32;
#     now it sits on G2:G37 viewing from the top ---
$ws.Range("G2:G37").Select()
$excel.ActiveWindow.ScrollRow = 1
